# Add a new slide at the end of the deck: "Prototyping" (Title and Content layout).
$p = $ppt.ActivePresentation

# Append a new slide using the same "Title and Content" layout (layout id 2)
# used by every other content slide in this deck.
$newIndex = $p.Slides.Count + 1
$slide = $p.Slides.Add($newIndex, 2)

# Title placeholder.
$title = $slide.Shapes.Item(1)
$title.Left = 53.33338742677165
$title.Top = 48.0
$title.Width = 676.9029921259843
$title.Height = 58.434804949606296
$title.TextFrame.TextRange.Text = "Prototyping"

# Content placeholder (left empty, matching the authored slide).
$content = $slide.Shapes.Item(2)
$content.Left = 53.33338742677165
$content.Top = 106.43488188976377
$content.Width = 676.9029921259843
$content.Height = 369.2629921259842
